$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.378699898719788
$ws.Range("B1").Value = 1.178440928459167
$ws.Range("C1").Value = 4.999112606048584
$ws.Range("D1").Value = 3.659501552581787
$ws.Range("E1").Value = 0.7109580039978027
